# This workbook is a weekly price-data extract. The update for this week
# inserts one new record at the top of the data block (row 19, right after
# the header row) and pushes every existing record down by one row, so the
# last existing record (previously row 123) becomes row 124 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 19; this shifts rows 19:123 down to 20:124
# (and keeps row 18's formatting to seed the new row's styles, matching the
# date/number formatting already used by column D throughout the table).
$ws.Rows("19").Insert()

# Populate the newly inserted row 19 with this week's new record.
$ws.Range("A19").Value = 5
$ws.Range("B19").Value = "Macroferia Regional de Talca"
$ws.Range("C19").Value = "Maule"
$ws.Range("D19").Value = 44550
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 15000
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 600
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
